$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Country column (C) with new test data for rows 3-6
$ws.Range("C3").Value = "US"
$ws.Range("C4").Value = "Germany"
$ws.Range("C5").Value = "Italy"
$ws.Range("C6").Value = "Belgium"

# Update selected cell to C6 (as reflected in the saved file's sheetView selection)
$ws.Range("C6").Select()
